$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Append rows 32-41 (subjects 31-40) following the existing pattern
# columns A-K are raw data, L-T mirror/derive them via formulas
# ------------------------------------------------------------------

$rowData = @(
  @{ row=32; A=31; B="teaching";   C="performing"; D="inst_t"; E="inst_p"; F="articulation"; G="dynamics";    H="stim_a"; I="stim_d" },
  @{ row=33; A=32; B="teaching";   C="performing"; D="inst_t"; E="inst_p"; F="dynamics";     G="articulation"; H="stim_d"; I="stim_a" },
  @{ row=34; A=33; B="performing"; C="teaching";   D="inst_p"; E="inst_t"; F="articulation"; G="dynamics";    H="stim_a"; I="stim_d" },
  @{ row=35; A=34; B="teaching";   C="performing"; D="inst_t"; E="inst_p"; F="articulation"; G="dynamics";    H="stim_a"; I="stim_d" },
  @{ row=36; A=35; B="teaching";   C="performing"; D="inst_t"; E="inst_p"; F="dynamics";     G="articulation"; H="stim_d"; I="stim_a" },
  @{ row=37; A=36; B="performing"; C="teaching";   D="inst_p"; E="inst_t"; F="articulation"; G="dynamics";    H="stim_a"; I="stim_d" },
  @{ row=38; A=37; B="teaching";   C="performing"; D="inst_t"; E="inst_p"; F="articulation"; G="dynamics";    H="stim_a"; I="stim_d" },
  @{ row=39; A=38; B="performing"; C="teaching";   D="inst_p"; E="inst_t"; F="articulation"; G="dynamics";    H="stim_a"; I="stim_d" },
  @{ row=40; A=39; B="teaching";   C="performing"; D="inst_t"; E="inst_p"; F="dynamics";     G="articulation"; H="stim_d"; I="stim_a" },
  @{ row=41; A=40; B="performing"; C="teaching";   D="inst_p"; E="inst_t"; F="dynamics";     G="articulation"; H="stim_d"; I="stim_a" }
)

foreach ($r in $rowData) {
  $n = $r.row

  $ws.Range("A$n").Value = $r.A
  $ws.Range("B$n").Value = $r.B
  $ws.Range("C$n").Value = $r.C
  $ws.Range("D$n").Value = $r.D
  $ws.Range("E$n").Value = $r.E
  $ws.Range("F$n").Value = $r.F
  $ws.Range("G$n").Value = $r.G
  $ws.Range("H$n").Value = $r.H
  $ws.Range("I$n").Value = $r.I

  # style (thin border) to match the other data rows
  $ws.Range("A$n`:I$n").Style = $ws.Range("A2:I2").Style

  $ws.Range("J$n").Value = ","
  $ws.Range("K$n").Value = ";"

  $ws.Range("L$n").Formula = "=CONCATENATE(A$n,J$n)"
  $ws.Range("M$n").Formula = "=B$n"
  $ws.Range("N$n").Formula = "=C$n"
  $ws.Range("O$n").Formula = "=D$n"
  $ws.Range("P$n").Formula = "=E$n"
  $ws.Range("Q$n").Formula = "=F$n"
  $ws.Range("R$n").Formula = "=G$n"
  $ws.Range("S$n").Formula = "=H$n"
  $ws.Range("T$n").Formula = "=CONCATENATE(I$n,K$n)"
}

# ------------------------------------------------------------------
# Update the print area to cover the newly added rows
# ------------------------------------------------------------------
$ws.PageSetup.PrintArea = '$A$1:$I$41'

# ------------------------------------------------------------------
# Update the visible selection to match the new extent
# ------------------------------------------------------------------
$ws.Range("A1:I41").Select()
